# [IMP] procurement result template
#
# The template's header row (row 6) had three English placeholder labels
# ("item" / "Asset Value" / "PO Total") that are replaced here with their
# Thai equivalents. In the shared-string table the three strings also swap
# positions (the old "item" slot now holds the text that used to be "PO
# Total", etc.) - but what matters for the workbook's visible content is
# simply which text ends up under which column, so we just assign the
# final Thai text directly to each header cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 header cells: I6 / J6 / K6
$ws.Range("I6").Value = "ครุภัณฑ์"
$ws.Range("J6").Value = "มูลค่าครุภัณฑ์"
$ws.Range("K6").Value = "ยอดรวมการจัดซื้อ"

# Reposition the outer application window (best effort - mirrors the
# author's workbookView change from the saved file metadata).
$win = $wb.Windows.Item(1)
$win.Left = 0
$win.Top = 460
$win.Width = 28800
$win.Height = 16540

# The author's last selection before saving moved to J13.
$ws.Range("J13").Select()
